# The "Sumgait" group of rows (A4:E6) holds three related records:
#   row 4: Russian original article   + its Russian-language URL
#   row 5: English translation title  + its English URL (unaffected)
#   row 6: US Congressional Record    + its congress.gov URL
#
# The source data (rebuilt from an upstream JSON feed per the commit
# message) reshuffled the shared-string table so that the text shown for
# row 4 and row 6 trade places, in both the title column (A) and the uri
# column (E), while row 5 stays untouched. Re-create that by swapping the
# cell values between row 4 and row 6 for columns A and E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$titleRow4 = $ws.Range("A4").Value2
$titleRow6 = $ws.Range("A6").Value2
$uriRow4   = $ws.Range("E4").Value2
$uriRow6   = $ws.Range("E6").Value2

$ws.Range("A4").Value = $titleRow6
$ws.Range("A6").Value = $titleRow4

$ws.Range("E4").Value = $uriRow6
$ws.Range("E6").Value = $uriRow4
